$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.171.72'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '1.826.40'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5983'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06945'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2746'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.29'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.02%  '
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").Value = '1.825.85'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.732'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6244'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009832'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.04%  '
$ws.Range("E16").Value = '  -5.11%  '
$ws.Range("D17").Value = '28.830.33'
$ws.Range("E17").Value = '  -1.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.554'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -11.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.871'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '156.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.919'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1282'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.86%  '
$ws.Range("E27").Value = '  -4.66%  '
$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06450'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.59%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.416'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.441'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.831'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.767'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.724'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.091'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6457'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.535'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.740'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01749'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.517'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.42%  '
$ws.Range("D40").Value = '1.142.50'
$ws.Range("E40").Value = '  -7.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8846'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").Value = '1.991.00'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.98'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000113'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.598'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.476'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.00%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4542'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05498'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.388'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.62%  '
